# Applies the "Added create text feauture" change to the API.xlsx workbook.
# This updates the DisplayState / ChangeStateName entries and appends three
# new function entries (DeleteState, SetTargetState, SetTargetStateObject)
# to the table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing cells
# ---------------------------------------------------------------------------

# DisplayState(int id) -> DisplayState()
$ws.Range("B17").Value = "DisplayState()"

# Example for DisplayState no longer passes the id argument
$ws.Range("B19").Value = "Example :SendMessage ('Facade', 'DisplayState');"
$ws.Range("B19").Characters(1, 7).Font.Bold = $true

# ChangeStateName description: "specified state" -> "target state"
$ws.Range("B21").Value = "Tells the editor to change the name of the target state. Requires you to pass a string and seprate arguments with a comma"

# ChangeStateName example updated to use a literal string argument
$ws.Range("B22").Value = "Example:SendMessage('Façade','ChangeStateName','this is a new name'); "
$ws.Range("B22").Characters(1, 7).Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. DeleteState() block (rows 23-25), formatted like the DisplayState block
# ---------------------------------------------------------------------------

$ws.Range("A17:B19").Copy() | Out-Null
$ws.Range("A23:B25").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:A25").Merge() | Out-Null

$ws.Range("A23").Value = "None"
$ws.Range("B23").Value = "DeleteState()"
$ws.Range("B24").Value = "Tells the editor to display to change the state to the specified id"
$ws.Range("B25").Value = "Example :SendMessage ('Facade', 'DeleteState');"
$ws.Range("B25").Characters(1, 7).Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. SetTargetState(int id) block (rows 26-28), formatted like the
#    ChangeStateName block
# ---------------------------------------------------------------------------

$ws.Range("A20:B22").Copy() | Out-Null
$ws.Range("A26:B28").PasteSpecial(-4122) | Out-Null
$ws.Range("A26:A28").Merge() | Out-Null

$ws.Range("A26").Value = "None"
$ws.Range("B26").Value = "SetTargetState(int id)"
$ws.Range("B27").Value = "Operations that target states would be targetting the specified state from now on"
$ws.Range("B28").Value = "Example :SendMessage ('Facade', SetTargetState',id);"
$ws.Range("B28").Characters(1, 9).Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. SetTargetStateObject(int id) block (rows 29-31)
# ---------------------------------------------------------------------------

$ws.Range("A17:B17").Copy() | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null

$ws.Range("A18").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null

$ws.Range("B19").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null

$ws.Range("A19:B19").Copy() | Out-Null
$ws.Range("A31:B31").PasteSpecial(-4122) | Out-Null

$ws.Range("A29:A31").Merge() | Out-Null

$ws.Range("A29").Value = "None"
$ws.Range("B29").Value = "SetTargetStateObject(int id)"
$ws.Range("B30").Value = "Operations that target states would be targetting the specified state object from now on"

# B31 repeats the exact same "Example" text as B19 (DisplayState) - copy the
# cell itself (not retype it) so the workbook reuses the same shared string.
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 5. Sheet bookkeeping: dimension / selected cell
# ---------------------------------------------------------------------------

$ws.Range("F24").Select() | Out-Null

Write-Host "edit applied"
